# The confusion-matrix figure had the "FP" and "FN" quadrant labels
# swapped (false positive / false negative were in the wrong boxes).
# Walk the shapes on the slide and swap the two mislabeled text runs.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $shp = $s.Shapes.Item($i)
    if ($shp.HasTextFrame) {
        $tr = $shp.TextFrame.TextRange
        if ($tr.Text -eq "FP") {
            $tr.Text = "FN"
        } elseif ($tr.Text -eq "FN") {
            $tr.Text = "FP"
        }
    }
}
